$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Add the new "messagePatient" row (row 4), following the existing pattern of
# role / email(hyperlink) / password columns used by rows 2-3.
$ws.Range("A4").Value = "messagePatient"
$ws.Range("B4").Value = "messagepatient@gmail.com"
$ws.Range("C4").Value = "LogBoxMaster"

# Turn the email into a mailto hyperlink, matching the style used on B2/B3.
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:messagepatient@gmail.com") | Out-Null

# The author subsequently stripped the default blue/underline hyperlink look
# from this particular cell while leaving the hyperlink itself intact.
$ws.Range("B4").Font.Underline = $false

# Reflect the final active selection on the sheet.
$ws.Range("B4").Select() | Out-Null

# Apply the page setup used for this worksheet (A4, portrait) as in the
# authored workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
